$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BasicSearch")

# Make BasicSearch the active/selected sheet (tab 1 becomes active tab)
$ws.Select()

# Update the search text cell
$ws.Range("B2").Value = "Apartments for Rent in Manhattan, NY"

# Move the sheet's selection to the edited cell
$ws.Range("B2").Select()

# Widen column B to fit the new (longer) text
$ws.Columns.Item(2).ColumnWidth = 42.5
